$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group and Self Assessment")

# Two previously-unassessed "placeholder" students (C14/C15) get their
# real student numbers filled in, replacing the generic "Student 5" /
# "Student 6" text placeholders.
$ws.Range("C14").Value2 = 1232250
$ws.Range("C15").Value2 = 1232233

# The team member on row 13 (student 1232225) fills in self/peer grades
# for the first six columns (D:I) with a grade of 5.
$ws.Range("D13:I13").Value2 = 5

# Reflect where the author was working when they saved.
$ws.Activate()
$ws.Range("G13").Select() | Out-Null
